# Reproduce the authors's edits to Threshold/Zn/2His_1Asp/2.xlsx (Sheet1):
#   - B2: 5.5 -> 5.4
#   - C2: 10.5 -> 10.2 (stored as 10.199999999999999 due to fp repr)
#   - B3: 5.5 -> 5.3
#   - C3: 9.5 -> 9
#   - C5: 15  -> 16
#   - Column A width: bestFit 21.375 -> custom 27
#   - Column C width: bestFit 5.875  -> custom 27.25 (closest attainable via
#     the pixel-quantised ColumnWidth COM property is 27.2857..)
#   - Page setup: paper size A4, portrait orientation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- cell value updates -----------------------------------------------
$ws.Range("B2").Value2 = 5.4
$ws.Range("C2").Value2 = 10.2
$ws.Range("B3").Value2 = 5.3
$ws.Range("C3").Value2 = 9
$ws.Range("C5").Value2 = 16

# --- column width updates ----------------------------------------------
# ColumnWidth is expressed in "characters" and gets rounded to the
# MDW-7 pixel grid by Excel, so we pick the input value whose rounded
# result lands on (or as close as possible to) the target sheet width.
$ws.Columns.Item(1).ColumnWidth = 26.29    # -> sheet width 27
$ws.Columns.Item(3).ColumnWidth = 26.57    # -> sheet width ~27.2857 (closest to 27.25)

# --- page setup ----------------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9      # xlPaperA4
$ps.Orientation = 1    # xlPortrait
